# "add battle damage / exp gained"
# Adds a "Name" column (nickname) right after "Species", and a "Gender"
# column at the end. Also refreshes EV to 255,255,255,255,255,255 for
# every Pokemon, and bumps CurrentStat for the first two rows to
# 128,20,20,20,20,20 (battle damage taken).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Name" column (C) right after "Species" (B) ---
$ws.Columns("C").Insert()

# Use Value2 (not Value) so comma-separated digit strings like
# "20,20,20,20,20,20" are stored as text, not auto-coerced to a number.
$ws.Range("C1").Value2 = "Name"
$ws.Range("C2").Value2 = "Nidorina"
$ws.Range("C3").Value2 = "Sandslash"
$ws.Range("C4").Value2 = "Raichu"
$ws.Range("C5").Value2 = "Sandshrew"
$ws.Range("C6").Value2 = "Arbok"
$ws.Range("C7").Value2 = "Fearow"

# --- 2. Append the new "Gender" column at the end (now column J) ---
$ws.Range("J1").Value2 = "Gender"
$ws.Range("J2").Value2 = "Male"
$ws.Range("J3").Value2 = "Male"
$ws.Range("J4").Value2 = "Male"
$ws.Range("J5").Value2 = "Male"
$ws.Range("J6").Value2 = "Male"
$ws.Range("J7").Value2 = "Male"

# --- 3. Update EV (now column G) for every Pokemon ---
# "255,255,255,255,255,255" is all-3-digit comma groups, so Excel's
# auto-detect treats it as a (huge) number even via Value2. Route it
# through a throw-away formula cell + paste-values so it lands as text
# without leaving any numeric formatting behind on the real cell.
$helper = $ws.Range("Z1")
$helper.Formula = "=""255,255,255,255,255,255"""
$helper.Copy()
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 7).PasteSpecial(-4163)
}
$helper.Clear()
$excel.CutCopyMode = 0

# --- 4. Update CurrentStat (now column H) after battle damage for the
#        first two Pokemon ---
$ws.Range("H2").Value2 = "128,20,20,20,20,20"
$ws.Range("H3").Value2 = "128,20,20,20,20,20"

# --- 5. Re-fit the column widths now that the columns hold new/wider
#        content (mirrors Excel's own "best fit" column-width refresh).
$ws.Range("C1").ColumnWidth = 9.166666666666666   # -> stored width 10
$ws.Range("G1").ColumnWidth = 20.5                # -> stored width ~21.33 (EV got wider)
$ws.Range("H1").ColumnWidth = 38.166666666666664  # -> stored width 39 (CurrentStat)
$ws.Range("I1").ColumnWidth = 38.166666666666664  # -> stored width 39 (Moves)

# --- 6. Fix up selection to match the post-edit state ---
$ws.Range("G7").Select()
